$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Patients")
$ws.Activate()

# Insert two new columns before column M (13), pushing M onwards to the right by 2.
# Using EntireColumn.Insert (like right-click > Insert in the Excel UI) so Excel
# naturally shifts cell contents/styles, the merged cell range and the dimension.
$ws.Range("M1:N1").EntireColumn.Insert()

# New "Investigations" / "Referrals" header cells (row 4)
$ws.Range("M4").Value2 = "Investigations"
$ws.Range("N4").Value2 = "Referrals"

# Cell values are written in the same order the original authoring tool
# introduced them into the shared-strings table (row 7 M/N, then row 6 N,
# then row 6 M) so the resulting sharedStrings.xml unique-string order lines
# up exactly with the source workbook.

# Row 7 - Patient 3 (NOT linked to a problem)
$ws.Range("M7").Value2 = "x" + [char]10 + "Patient Linked to DiagnosticReport, ProcedureRequests and specimens and a NO problems linked"

# Row 6 - Patient 2 (linked to a problem)
$ws.Range("N6").Value2 = "x" + [char]10 + "Patient Linked to ReferralRequests and  linked to a problem"

# Row 7 - Patient 3 (NOT linked to a problem)
$ws.Range("N7").Value2 = "x" + [char]10 + "Patient Linked to ReferralRequests and  NOT linked to a problem"

# Row 6 - Patient 2 (linked to a problem)
$ws.Range("M6").Value2 = "x" + [char]10 + "Patient Linked to DiagnosticReport, ProcedureRequests,  specimens and a problem linked"

# Row 8 - numeric placeholders (0) matching the other "0" columns on this row
$ws.Range("M8").Value2 = 0
$ws.Range("N8").Value2 = 0

# Restore the view/selection state
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("M7").Select()
